$d = $word.ActiveDocument

$pairs = @(
    @("338×3=", "625×5="),
    @("807×6=", "815×8="),
    @("574×7=", "223×3="),
    @("783×5=", "590×3="),
    @("530×3=", "508×2="),
    @("660×4=", "393×9="),
    @("682×8=", "756×5="),
    @("843×6=", "566×4="),
    @("972×2=", "404×7="),
    @("336×8=", "456×5="),
    @("771×6=", "983×9="),
    @("511×3=", "553×8="),
    @("368×9=", "828×5="),
    @("942×2=", "791×5="),
    @("740×5=", "942×6="),
    @("808×3=", "273×2="),
    @("124×6=", "500×2="),
    @("788×6=", "647×9="),
    @("398×8=", "263×4="),
    @("996×5=", "803×8="),
    @("215×6=", "865×7="),
    @("263×5=", "679×7="),
    @("154×5=", "874×5="),
    @("307×7=", "638×4="),
    @("738×6=", "449×9=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
